$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "jvjvj"
$ws.Range("A4").Value = "bjjh"
$ws.Range("A4").Select()
